$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (col A / col C get wider, to fit the new long test-step text) ---
$ws.Columns.Item(1).ColumnWidth = 59.17
$ws.Columns.Item(3).ColumnWidth = 49.17

# --- Remove the old per-page detail rows (rows 6-14), shifting everything up ---
$ws.Rows("6:14").Delete() | Out-Null

# --- Row 1: add xpath-selector header cells in C1:G1 (one per form section) ---
$ws.Range("C1").Value = '//*[@id="insurance-form"]/div/section[1]'
$ws.Range("D1").Value = '//*[@id="insurance-form"]/div/section[2]'
$ws.Range("E1").Value = '//*[@id="insurance-form"]/div/section[3]'
$ws.Range("F1").Value = '//*[@id="insurance-form"]/div/section[4]'
$ws.Range("G1").Value = '//*[@id="insurance-form"]/div/section[5]'

# --- Row 2: highlight C2:G2 with a yellow fill (new, empty styled cells) ---
$ws.Range("C2:G2").Interior.Color = 65535

# --- Row 4 ("Check defaults"): add SELECTED / NOTSELECTED flags across C4:G4 ---
$ws.Range("C4").Value = "<SELECTED>"
$ws.Range("D4").Value = "<NOTSELECTED>"
$ws.Range("E4").Value = "<NOTSELECTED>"
$ws.Range("F4").Value = "<NOTSELECTED>"
$ws.Range("G4").Value = "<NOTSELECTED>"

# --- Row 5: replace the old "Auswahl Page Vehicle Data" row with the new smoke-test row ---
$ws.Range("A5").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPageVehicleData"
$ws.Range("C5").Value = "102_VehicleInsuranceAutomobile_001_SmokeTest_FillPage"

# --- Reposition the screenshot picture so it again starts a couple of rows below the
#     (now much shorter) table, preserving its original on-screen size ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 100.80007874015747
$shp.Left = 0
$shp.Width = 676.6974015748032
$shp.Height = 398.4

# --- Selection, as left by the author ---
$ws.Range("F20").Select() | Out-Null
